$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers to Excel's auto-detection
# (e.g. "609.87", "1.00") need to be forced to Text format first so the
# literal string (including trailing zeros) is preserved exactly, matching
# the source report data instead of being parsed into a numeric value.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D28", "D29", "D30", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D41", "D43", "D44", "D45", "D47", "D48", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply each cell's new value
$ws.Range('D2').Value = '66.378.95'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '3.248.59'
$ws.Range('E3').Value = '  +2.54%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '609.87'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').Value = '156.89'
$ws.Range('E6').Value = '  +1.89%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.247.04'
$ws.Range('E8').Value = '  +2.61%  '
$ws.Range('D9').Value = '0.544'
$ws.Range('E9').Value = '  -0.34%  '
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').Value = '  +1.92%  '
$ws.Range('D11').Value = '5.78'
$ws.Range('E11').Value = '  +3.38%  '
$ws.Range('D12').Value = '0.496'
$ws.Range('E12').Value = '  -3.97%  '
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('D14').Value = '38.97'
$ws.Range('E14').Value = '  +1.84%  '
$ws.Range('D15').Value = '3.784.62'
$ws.Range('E15').Value = '  +2.63%  '
$ws.Range('D16').Value = '66.500.37'
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('D17').Value = '7.44'
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D18').Value = '3.248.44'
$ws.Range('E18').Value = '  +2.56%  '
$ws.Range('E19').Value = '  +1.09%  '
$ws.Range('D20').Value = '504.52'
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('D21').Value = '15.42'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').Value = '0.750'
$ws.Range('E22').Value = '  +3.25%  '
$ws.Range('D23').Value = '8.08'
$ws.Range('E23').Value = '  +1.14%  '
$ws.Range('D24').Value = '14.69'
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('D25').Value = '87.05'
$ws.Range('E25').Value = '  +2.88%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  +0.56%  '
$ws.Range('D28').Value = '9.16'
$ws.Range('E28').Value = '  +0.84%  '
$ws.Range('D29').Value = '2.38'
$ws.Range('E29').Value = '  +0.30%  '
$ws.Range('D30').Value = '0.127'
$ws.Range('E30').Value = '  +44.28%  '
$ws.Range('D31').Value = '6.95'
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('D32').Value = '2.88'
$ws.Range('E32').Value = '  -5.02%  '
$ws.Range('D33').Value = '27.89'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('E35').Value = '  -3.37%  '
$ws.Range('D36').Value = '6.44'
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('D37').Value = '55.47'
$ws.Range('E37').Value = '  +1.06%  '
$ws.Range('D38').Value = '3.32'
$ws.Range('E38').Value = '  +17.92%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = '495.20'
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0778'
$ws.Range('E40').Value = '  +13.91%  '
$ws.Range('D41').Value = '0.0422'
$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('E42').Value = '  +1.20%  '
$ws.Range('D43').Value = '8.82'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('D44').Value = '2.51'
$ws.Range('E44').Value = '  +3.77%  '
$ws.Range('D45').Value = '0.292'
$ws.Range('E45').Value = '  -1.37%  '
$ws.Range('D46').Value = '2.978.35'
$ws.Range('E46').Value = '  +5.45%  '
$ws.Range('D47').Value = '28.83'
$ws.Range('E47').Value = '  +3.22%  '
$ws.Range('D48').Value = '2.50'
$ws.Range('E48').Value = '  +5.68%  '
$ws.Range('E49').Value = '  +2.20%  '
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '120.82'
$ws.Range('E51').Value = '  -0.95%  '
